$d = $word.ActiveDocument

# --- Edit 1: insert a new empty "Body Text" paragraph right before the
# "Class meetings: Online, asynchronous ..." paragraph (splits what was a
# single "Description..." paragraph boundary so a blank Body Text
# paragraph now precedes "Class meetings:"). We find the paragraph that
# ends right before "Class meetings:" and insert a paragraph break at its
# end, which inherits the BodyText style.
$classMeetingsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Class meetings:*Online, asynchronous*") {
        $classMeetingsPara = $p
        break
    }
}
$prevEnd = $classMeetingsPara.Previous().Range.End
$insertRange = $d.Range($prevEnd, $prevEnd)
$insertRange.InsertAfter("`r")

# --- Edit 2: remove the empty "First Paragraph" styled paragraph that
# sits between "develop step-by-step problem solving and debugging
# practices" and the "Required Text" heading.
$emptyPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Previous() -ne $null -and $p.Next() -ne $null) {
        if (($p.Previous().Range.Text -like "*develop step-by-step problem solving and debugging practices*") -and ($p.Next().Range.Text -like "Required Text*")) {
            $emptyPara = $p
            break
        }
    }
}
$emptyPara.Range.Delete()
